# Insert a new data row at sheet row 55, shifting existing rows 55..195 down to 56..196.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("55").Insert()

# Fill in the newly inserted row 55 with the new weekly record.
$ws.Range("A55").Value = 7
$ws.Range("B55").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C55").Value = "Ñuble"
$ws.Range("D55").Value = 44498
$ws.Range("E55").Value = 16
$ws.Range("F55").Value = 100114013
$ws.Range("G55").Value = "Zanahoria"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 120
$ws.Range("K55").Value = 7500
$ws.Range("L55").Value = 8000
$ws.Range("M55").Value = 7750
$ws.Range("N55").Value = '$/saco 20 kilos'
$ws.Range("O55").Value = "Provincia de Diguillín"
$ws.Range("P55").Value = 388
$ws.Range("Q55").Value = 20
$ws.Range("R55").Value = "Hortaliza"
